$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G3").Value = 4.33
$ws.Range("H3").Value = 3
$ws.Range("I3").Value = 2
$ws.Range("K3").Value = 1.91
$ws.Range("W3").Value = 9
$ws.Range("X3").Value = 19
$ws.Range("Y3").Value = 15
$ws.Range("Z3").Value = 41
$ws.Range("AC3").Value = 6
$ws.Range("AK3").Value = 21
$ws.Range("AX3").Value = 12
$ws.Range("BA3").Value = 81
